$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Header row (row 1): turn it into a full header label row B..M -----
# (the former row 1 held the first data row's own values; the new layout
# reuses row 1 purely for column headers, and the real values move to row 2)
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Give the new header cells (F1:M1) the same look (bold font, border,
# centered) as the pre-existing header cells by copying the format from E1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:M1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ----- Data row (row 2): keep existing values, add the new columns -----
$ws.Range("B2").Value = "台新國際商業銀行南京東路分行"
$ws.Range("C2").Value = "活期存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "李慶華"
$ws.Range("F2").Value = 7175

$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
# Leading apostrophe forces the date-looking string to be kept as literal text
# instead of being auto-converted into a date serial number. Re-apply the
# "Normal" style afterwards so the cell doesn't end up with a stray style
# index just for the quote-prefix marker.
$ws.Range("I2").Value = "'2011-11-30"
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = "李慶華"
$ws.Range("K2").Value = 607
$ws.Range("L2").Value = "tmp2e001"
$ws.Range("M2").Value = 12
